$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).NumberFormat = '@'
$ws.Cells.Item(2, 4).Value = '65.305.84'
$ws.Cells.Item(2, 5).Value = '  -1.66%  '

$ws.Cells.Item(3, 4).NumberFormat = '@'
$ws.Cells.Item(3, 4).Value = '3.424.40'
$ws.Cells.Item(3, 5).Value = '  -4.77%  '

$ws.Cells.Item(4, 5).Value = '  +0.03%  '

$ws.Cells.Item(5, 4).NumberFormat = '@'
$ws.Cells.Item(5, 4).Value = '593.86'
$ws.Cells.Item(5, 5).Value = '  -2.26%  '

$ws.Cells.Item(6, 4).NumberFormat = '@'
$ws.Cells.Item(6, 4).Value = '134.00'
$ws.Cells.Item(6, 5).Value = '  -9.70%  '

$ws.Cells.Item(7, 4).NumberFormat = '@'
$ws.Cells.Item(7, 4).Value = '3.422.01'
$ws.Cells.Item(7, 5).Value = '  -4.85%  '

$ws.Cells.Item(8, 4).NumberFormat = '@'
$ws.Cells.Item(8, 4).Value = '0.999'
$ws.Cells.Item(8, 5).Value = '  -0.12%  '

$ws.Cells.Item(9, 4).NumberFormat = '@'
$ws.Cells.Item(9, 4).Value = '0.488'
$ws.Cells.Item(9, 5).Value = '  +0.20%  '

$ws.Cells.Item(10, 4).NumberFormat = '@'
$ws.Cells.Item(10, 4).Value = '7.46'
$ws.Cells.Item(10, 5).Value = '  -5.11%  '

$ws.Cells.Item(11, 5).Value = '  -10.38%  '

$ws.Cells.Item(12, 4).NumberFormat = '@'
$ws.Cells.Item(12, 4).Value = '0.376'
$ws.Cells.Item(12, 5).Value = '  -9.15%  '

$ws.Cells.Item(13, 4).NumberFormat = '@'
$ws.Cells.Item(13, 4).Value = '4.001.32'
$ws.Cells.Item(13, 5).Value = '  -4.84%  '

$ws.Cells.Item(14, 4).NumberFormat = '@'
$ws.Cells.Item(14, 4).Value = '0.0000179'
$ws.Cells.Item(14, 5).Value = '  -13.10%  '

$ws.Cells.Item(15, 4).NumberFormat = '@'
$ws.Cells.Item(15, 4).Value = '26.28'
$ws.Cells.Item(15, 5).Value = '  -11.11%  '

$ws.Cells.Item(16, 2).Value = 'WrappedBTC'
$ws.Cells.Item(16, 3).Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Cells.Item(16, 4).NumberFormat = '@'
$ws.Cells.Item(16, 4).Value = '65.299.45'
$ws.Cells.Item(16, 5).Value = '  -1.72%  '

$ws.Cells.Item(17, 2).Value = 'WrappedEther'
$ws.Cells.Item(17, 3).Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Cells.Item(17, 4).NumberFormat = '@'
$ws.Cells.Item(17, 4).Value = '3.413.61'
$ws.Cells.Item(17, 5).Value = '  -4.97%  '

$ws.Cells.Item(18, 5).Value = '  -3.02%  '

$ws.Cells.Item(19, 4).NumberFormat = '@'
$ws.Cells.Item(19, 4).Value = '9.89'
$ws.Cells.Item(19, 5).Value = '  -10.75%  '

$ws.Cells.Item(20, 4).NumberFormat = '@'
$ws.Cells.Item(20, 4).Value = '5.71'
$ws.Cells.Item(20, 5).Value = '  -9.74%  '

$ws.Cells.Item(21, 4).NumberFormat = '@'
$ws.Cells.Item(21, 4).Value = '13.61'
$ws.Cells.Item(21, 5).Value = '  -8.56%  '

$ws.Cells.Item(22, 4).NumberFormat = '@'
$ws.Cells.Item(22, 4).Value = '391.62'
$ws.Cells.Item(22, 5).Value = '  -7.39%  '

$ws.Cells.Item(23, 4).NumberFormat = '@'
$ws.Cells.Item(23, 4).Value = '73.07'
$ws.Cells.Item(23, 5).Value = '  -7.13%  '

$ws.Cells.Item(24, 4).NumberFormat = '@'
$ws.Cells.Item(24, 4).Value = '0.541'
$ws.Cells.Item(24, 5).Value = '  -11.50%  '

$ws.Cells.Item(25, 5).Value = '  -0.20%  '

$ws.Cells.Item(26, 4).NumberFormat = '@'
$ws.Cells.Item(26, 4).Value = '3.567.15'
$ws.Cells.Item(26, 5).Value = '  -4.63%  '

$ws.Cells.Item(27, 5).Value = '  -13.27%  '

$ws.Cells.Item(28, 4).NumberFormat = '@'
$ws.Cells.Item(28, 4).Value = '1.00'
$ws.Cells.Item(28, 5).Value = '  -0.01%  '

$ws.Cells.Item(29, 5).Value = '  -10.04%  '

$ws.Cells.Item(30, 4).NumberFormat = '@'
$ws.Cells.Item(30, 4).Value = '7.10'
$ws.Cells.Item(30, 5).Value = '  -13.85%  '

$ws.Cells.Item(31, 4).NumberFormat = '@'
$ws.Cells.Item(31, 4).Value = '8.14'
$ws.Cells.Item(31, 5).Value = '  -13.28%  '

$ws.Cells.Item(32, 4).NumberFormat = '@'
$ws.Cells.Item(32, 4).Value = '3.430.69'
$ws.Cells.Item(32, 5).Value = '  -4.51%  '

$ws.Cells.Item(33, 5).Value = '  -0.03%  '

$ws.Cells.Item(34, 5).Value = '  -8.15%  '

$ws.Cells.Item(35, 4).NumberFormat = '@'
$ws.Cells.Item(35, 4).Value = '22.58'
$ws.Cells.Item(35, 5).Value = '  -10.07%  '

$ws.Cells.Item(36, 2).Value = 'Monero'
$ws.Cells.Item(36, 3).Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Cells.Item(36, 4).NumberFormat = '@'
$ws.Cells.Item(36, 4).Value = '171.66'
$ws.Cells.Item(36, 5).Value = '  -1.95%  '

$ws.Cells.Item(37, 2).Value = 'Fetch.AI'
$ws.Cells.Item(37, 3).Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Cells.Item(37, 4).NumberFormat = '@'
$ws.Cells.Item(37, 4).Value = '1.22'
$ws.Cells.Item(37, 5).Value = '  -15.14%  '

$ws.Cells.Item(38, 4).NumberFormat = '@'
$ws.Cells.Item(38, 4).Value = '6.77'
$ws.Cells.Item(38, 5).Value = '  -12.81%  '

$ws.Cells.Item(39, 4).NumberFormat = '@'
$ws.Cells.Item(39, 4).Value = '1.51'
$ws.Cells.Item(39, 5).Value = '  -9.75%  '

$ws.Cells.Item(40, 4).NumberFormat = '@'
$ws.Cells.Item(40, 4).Value = '4.80'
$ws.Cells.Item(40, 5).Value = '  -14.09%  '

$ws.Cells.Item(41, 4).NumberFormat = '@'
$ws.Cells.Item(41, 4).Value = '0.0766'
$ws.Cells.Item(41, 5).Value = '  -10.14%  '

$ws.Cells.Item(42, 4).NumberFormat = '@'
$ws.Cells.Item(42, 4).Value = '0.809'
$ws.Cells.Item(42, 5).Value = '  -8.51%  '

$ws.Cells.Item(43, 4).NumberFormat = '@'
$ws.Cells.Item(43, 4).Value = '43.45'
$ws.Cells.Item(43, 5).Value = '  -5.47%  '

$ws.Cells.Item(44, 4).NumberFormat = '@'
$ws.Cells.Item(44, 4).Value = '1.00'
$ws.Cells.Item(44, 5).Value = '  +0.21%  '

$ws.Cells.Item(45, 4).NumberFormat = '@'
$ws.Cells.Item(45, 4).Value = '4.37'
$ws.Cells.Item(45, 5).Value = '  -15.75%  '

$ws.Cells.Item(46, 4).NumberFormat = '@'
$ws.Cells.Item(46, 4).Value = '1.60'
$ws.Cells.Item(46, 5).Value = '  -13.57%  '

$ws.Cells.Item(47, 4).NumberFormat = '@'
$ws.Cells.Item(47, 4).Value = '1.08'
$ws.Cells.Item(47, 5).Value = '  -4.63%  '

$ws.Cells.Item(48, 4).NumberFormat = '@'
$ws.Cells.Item(48, 4).Value = '21.85'
$ws.Cells.Item(48, 5).Value = '  -8.25%  '

$ws.Cells.Item(49, 4).NumberFormat = '@'
$ws.Cells.Item(49, 4).Value = '6.50'
$ws.Cells.Item(49, 5).Value = '  -9.07%  '

$ws.Cells.Item(50, 4).NumberFormat = '@'
$ws.Cells.Item(50, 4).Value = '2.12'
$ws.Cells.Item(50, 5).Value = '  -16.30%  '

$ws.Cells.Item(51, 4).NumberFormat = '@'
$ws.Cells.Item(51, 4).Value = '2.187.11'
$ws.Cells.Item(51, 5).Value = '  -8.54%  '
